$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.820666790008545
$ws.Range("B1").Value = 2.798565864562988
$ws.Range("C1").Value = 2.016470432281494
$ws.Range("D1").Value = 1.859195709228516
$ws.Range("E1").Value = 1.915044069290161
